# Remove the "organisation" column from the template.
# Original header order: A=organisation, B=code, C=name, D=description, E=budget_usd
# Target header order:   A=code, B=name, C=description, D=budget_usd

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the remaining header values one column to the left, dropping "organisation".
$ws.Range("A1").Value = $ws.Range("B1").Value()
$ws.Range("B1").Value = $ws.Range("C1").Value()
$ws.Range("C1").Value = $ws.Range("D1").Value()
$ws.Range("D1").Value = $ws.Range("E1").Value()

# Remove the now-duplicate trailing column (this also drops the unused
# "organisation" shared string from the workbook).
$ws.Columns.Item(5).Delete()

# Match the resulting selection left by the edit (whole of column A selected).
$ws.Range("A1:A1048576").Select()
